$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-07 (row 20)
$ws.Range("B20").Value = 6199
$ws.Range("D20").Value = 5581316
$ws.Range("E20").Value = 900.3574770124213
$ws.Range("F20").Value = 7.082397650716876
$ws.Range("H20").Value = 26.25489376843753
